$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.859.44"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "1.630.13"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +0.67%  "
$ws.Range("D5").Value = "'214.77"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("D8").Value = "'28.70"
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'0.0608"
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("D12").Value = "1.862.72"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").Value = "1.635.50"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("D14").Value = "'0.568"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "'9.41"
$ws.Range("D16").Value = "29.873.11"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("D17").Value = "'3.83"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "'65.25"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("D19").Value = "'240.56"
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").Value = "0.0₃0701"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'4.13"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").Value = "'9.79"
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("D24").Value = "'2.18"
$ws.Range("E24").Value = "  +2.86%  "
$ws.Range("D25").Value = "'157.55"
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("D26").Value = "'15.48"
$ws.Range("E26").Value = "  -0.94%  "
$ws.Range("E27").Value = "  -1.01%  "
$ws.Range("D28").Value = "'6.60"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("D31").Value = "'1.11"
$ws.Range("E31").Value = "  +2.61%  "
$ws.Range("E32").Value = "  +1.97%  "
$ws.Range("E33").Value = "  -0.86%  "
$ws.Range("D34").Value = "1.426.05"
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("E35").Value = "  +3.30%  "
$ws.Range("E36").Value = "  -2.52%  "
$ws.Range("E37").Value = "  -3.91%  "
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("D40").Value = "'0.555"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("D41").Value = "'74.88"
$ws.Range("E41").Value = "  +7.00%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.0502"
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'0.831"
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("E46").Value = "  +0.79%  "
$ws.Range("D47").Value = "1.771.04"
$ws.Range("E48").Value = "  -2.10%  "
$ws.Range("D49").Value = "'48.84"
$ws.Range("E49").Value = "  -8.31%  "
$ws.Range("D50").Value = "'92.08"
$ws.Range("E50").Value = "  +4.74%  "
$ws.Range("E51").Value = "  +4.46%  "
